$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Harmonize the "Tags" section so that accession numbers use the short
# CURIE form (PREFIX:ID) stored only in the "Tags Term Accession Number"
# row, matching how the other tags are already formatted. This clears the
# now-redundant "Tags Term Source REF" entries for those tags.

# Tag 4 ("Mass Spectrometry" stays the same, NCIT_C17156 url -> NCIT:C17156 CURIE)
$ws.Range("D13").Value = "NCIT:C17156"
$ws.Range("D14").Value = ""

# Tag 3 ("assay protocol" -> "assay", DPBO:1000177 url -> OBI:0000070 CURIE)
$ws.Range("C12").Value = "assay"
$ws.Range("C13").Value = "OBI:0000070"
$ws.Range("C14").Value = ""
